# Database inputs feedback corrected and commented code deleted
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hydro")

# Correct the feedback database inputs for the Hydro turbine (D6, D7).
# Downstream formulas (D11-D14, D16) depend on these and recalculate automatically.
$ws.Range("D6").Value = 500
$ws.Range("D7").Value = 31

# Move the active selection to D7, matching the reviewed/edited cell.
$ws.Activate()
$ws.Range("D7").Select() | Out-Null
